# Forgot password and otp api Basic added
# Inserts 3 new rows (User: Change Password / Forgot Password OTP / Forgot
# Password validate OTP) above the existing row 12 ("Get All Vehicle Types
# From Master"), pushing the rest of the table down by 3 rows, then fixes
# up the AutoFilter range, the hidden _FilterDatabase defined name, and the
# active selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 blank rows before row 12 (data below shifts down to 15..28).
$ws.Rows("12:14").Insert()

# 2. Copy the formatting (styles) of the row above (row 11) down into the
#    three new rows so they match the rest of the table visually.
$ws.Range("B11:O11").Copy()
$ws.Range("B12:O14").PasteSpecial(-4122) # xlPasteFormats

# 3. Populate the three new data rows.

# Row 12: User / Change Password
$ws.Cells.Item(12,2).Value2  = "User"
$ws.Cells.Item(12,3).Value2  = "Change Password"
$ws.Cells.Item(12,4).Value2  = "WS-UP-07"
$ws.Cells.Item(12,5).Value2  = "app.user.password.update"
$ws.Cells.Item(12,6).Value2  = $true
$ws.Cells.Item(12,7).Value2  = "user"
$ws.Cells.Item(12,8).Value2  = "/passwordupdate"
$ws.Cells.Item(12,9).Value2  = "PUT"
$ws.Cells.Item(12,10).Value2 = "P1"
$ws.Cells.Item(12,11).Value2 = "Basic Done"

# Row 13: User / Forgot Password OTP
$ws.Cells.Item(13,2).Value2  = "User"
$ws.Cells.Item(13,3).Value2  = "Forgot Password OTP"
$ws.Cells.Item(13,4).Value2  = "WS-UP-08"
$ws.Cells.Item(13,5).Value2  = "app.user.otp.save"
$ws.Cells.Item(13,6).Value2  = $true
$ws.Cells.Item(13,7).Value2  = "user"
$ws.Cells.Item(13,8).Value2  = "/sendotp"
$ws.Cells.Item(13,9).Value2  = "POST"
$ws.Cells.Item(13,10).Value2 = "P1"
$ws.Cells.Item(13,11).Value2 = "Basic Done"

# Row 14: User / Forgot Password validate OTP
$ws.Cells.Item(14,2).Value2  = "User"
$ws.Cells.Item(14,3).Value2  = "Forgot Password validate OTP"
$ws.Cells.Item(14,4).Value2  = "WS-UP-09"
$ws.Cells.Item(14,5).Value2  = "app.user.otp.get"
$ws.Cells.Item(14,6).Value2  = $false
$ws.Cells.Item(14,7).Value2  = "user"
$ws.Cells.Item(14,8).Value2  = "/validateotp"
$ws.Cells.Item(14,9).Value2  = "POST"
$ws.Cells.Item(14,10).Value2 = "P1"
$ws.Cells.Item(14,11).Value2 = "Endpoint Added"

# 4. Re-create the generated-code formulas (columns M/N/O) for the 3 new
#    rows, matching the pattern used throughout the sheet.
foreach ($r in 12..14) {
    $ws.Cells.Item($r,13).Formula = '=_xlfn.CONCAT("INSERT INTO ",CHAR(34),"M_CTL_CONFIG",CHAR(34)," VALUES(''",D' + $r + ',"'',''CONNON_CONFIG'', 0, ''",C' + $r + ',"'', ''{}'', 0, 0, CURRENT_TIMESTAMP, ''ATUL'', null, null);")'
    $ws.Cells.Item($r,14).Formula = '=_xlfn.CONCAT(IF(I' + $r + '="GET","@GetMapping(",IF(I' + $r + '="POST","@PostMapping(",IF(I' + $r + '="DELETE","@DeleteMapping(",IF(I' + $r + '="PUT","@PutMapping(","")))),CHAR(34),H' + $r + ',CHAR(34),")")'
    $ws.Cells.Item($r,15).Formula = '=_xlfn.CONCAT("@ServiceInfo(serviceCode = ",CHAR(34),D' + $r + ',,CHAR(34),", serviceName = ",CHAR(34),C' + $r + ',CHAR(34), ", queryId = ",CHAR(34),E' + $r + ',CHAR(34),", logActivity =",F' + $r + ',")")'
}

# 5. Re-apply the AutoFilter over the new, larger range (B3:L28).
$ws.AutoFilterMode = $false
$ws.Range("B3:L28").AutoFilter()

# 6. Fix the hidden _xlnm._FilterDatabase defined name to match.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=ServicesList!`$B`$3:`$L`$28"
    }
}

# 7. Move the active selection to K15 (matches the edited file).
$ws.Range("K15").Select()
